$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 249
$ws.Range("I6").Value = 249
$ws.Range("K6").Value = 747
$ws.Range("M6").Value = -635
$ws.Range("H8").Value = 3.857143
$ws.Range("I8").Value = 3.857143
$ws.Range("K8").Value = 11.571429
$ws.Range("M8").Value = 127.428571
$ws.Range("H18").Value = 111112650
$ws.Range("I18").Value = 1468.3334
$ws.Range("K18").Value = 1468.3334
$ws.Range("M18").Value = -1184.3334
$ws.Range("H46").Value = 500
$ws.Range("I46").Value = 500
$ws.Range("J46").Value = 500
$ws.Range("K46").Value = 1500
$ws.Range("L46").Value = 1500
$ws.Range("M46").Value = -1381
$ws.Range("N46").Value = -1738
$ws.Range("H60").Value = 500
$ws.Range("I60").Value = 500
$ws.Range("J60").Value = 500
$ws.Range("K60").Value = 1500
$ws.Range("L60").Value = 1500
$ws.Range("M60").Value = -1016
$ws.Range("N60").Value = -2468
$ws.Range("H98").Value = 2607.261
$ws.Range("I98").Value = 2707.7144
$ws.Range("K98").Value = 2707.7144
$ws.Range("M98").Value = -1209.7144
$ws.Range("H111").Value = 5279.857
$ws.Range("I111").Value = 5909.8335
$ws.Range("J111").Value = 1500
$ws.Range("K111").Value = 17729.5005
$ws.Range("L111").Value = 4500
$ws.Range("M111").Value = -14662.5005
$ws.Range("N111").Value = -10634
$ws.Range("H112").Value = 1651.4762
$ws.Range("J112").Value = 1683.2632
$ws.Range("L112").Value = 5049.7896
$ws.Range("N112").Value = -7265.7896
$ws.Range("H122").Value = 2607.261
$ws.Range("I122").Value = 2707.7144
$ws.Range("K122").Value = 8123.1432
$ws.Range("M122").Value = -5673.1432
$ws.Range("H131").Value = 1253672.2
$ws.Range("J131").Value = 3879.6
$ws.Range("L131").Value = 11638.8
$ws.Range("N131").Value = -21718.8
$ws.Range("H132").Value = 4654.879
$ws.Range("I132").Value = 4769.125
$ws.Range("J132").Value = 999
$ws.Range("K132").Value = 14307.375
$ws.Range("L132").Value = 2997
$ws.Range("M132").Value = -11777.375
$ws.Range("N132").Value = -8057
$ws.Range("H135").Value = 2076.6667
$ws.Range("I135").Value = 1269
$ws.Range("K135").Value = 11421
$ws.Range("M135").Value = -8886
$ws.Range("H137").Value = 2051.25
$ws.Range("I137").Value = 1926.5
$ws.Range("J137").Value = 2924.5
$ws.Range("K137").Value = 5779.5
$ws.Range("L137").Value = 8773.5
$ws.Range("M137").Value = -3229.5
$ws.Range("N137").Value = -13873.5
$ws.Range("H141").Value = 3178.7693
$ws.Range("I141").Value = 2822
$ws.Range("J141").Value = 4368
$ws.Range("K141").Value = 8466
$ws.Range("L141").Value = 13104
$ws.Range("M141").Value = -3286
$ws.Range("N141").Value = -23464

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2499
$ws.Range("I2").Value = 2499
$ws.Range("K2").Value = 2499
$ws.Range("M2").Value = -2386
$ws.Range("H32").Value = 1661692.6
$ws.Range("I32").Value = 807032.3
$ws.Range("K32").Value = 807032.3
$ws.Range("M32").Value = -806745.3
$ws.Range("H45").Value = 29545.533
$ws.Range("I45").Value = 36481.168
$ws.Range("K45").Value = 36481.168
$ws.Range("M45").Value = -36104.168
$ws.Range("H74").Value = 1617.9375
$ws.Range("I74").Value = 853.36365
$ws.Range("K74").Value = 853.36365
$ws.Range("M74").Value = 20.63634999999999
$ws.Range("H77").Value = 1617.9375
$ws.Range("I77").Value = 853.36365
$ws.Range("K77").Value = 4266.81825
$ws.Range("M77").Value = 101.1817499999997
$ws.Range("H80").Value = 69403.336
$ws.Range("J80").Value = 94055
$ws.Range("L80").Value = 94055
$ws.Range("N80").Value = -96051
$ws.Range("H83").Value = 69403.336
$ws.Range("J83").Value = 94055
$ws.Range("L83").Value = 282165
$ws.Range("N83").Value = -292149
$ws.Range("H110").Value = 614
$ws.Range("I110").Value = 522.4286
$ws.Range("K110").Value = 522.4286
$ws.Range("M110").Value = 1522.5714
$ws.Range("H116").Value = 2499
$ws.Range("I116").Value = 2499
$ws.Range("K116").Value = 2499
$ws.Range("M116").Value = -205

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2499
$ws.Range("I3").Value = 2499
$ws.Range("K3").Value = 2499
$ws.Range("M3").Value = -2385
$ws.Range("H86").Value = 2936.3823
$ws.Range("I86").Value = 3085.1904
$ws.Range("J86").Value = 2696
$ws.Range("K86").Value = 3085.1904
$ws.Range("L86").Value = 2696
$ws.Range("M86").Value = -1962.1904
$ws.Range("N86").Value = -4942
$ws.Range("H89").Value = 2936.3823
$ws.Range("I89").Value = 3085.1904
$ws.Range("J89").Value = 2696
$ws.Range("K89").Value = 15425.952
$ws.Range("L89").Value = 13480
$ws.Range("M89").Value = -9809.951999999999
$ws.Range("N89").Value = -24712
$ws.Range("H134").Value = 2078.32
$ws.Range("I134").Value = 1439.25
$ws.Range("J134").Value = 2668.2307
$ws.Range("K134").Value = 4317.75
$ws.Range("L134").Value = 8004.6921
$ws.Range("M134").Value = -1782.75
$ws.Range("N134").Value = -13074.6921

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 235.35294
$ws.Range("I7").Value = 253.36363
$ws.Range("K7").Value = 253.36363
$ws.Range("M7").Value = -140.36363
$ws.Range("H16").Value = 748.0833
$ws.Range("I16").Value = 748.0833
$ws.Range("K16").Value = 748.0833
$ws.Range("M16").Value = -461.0833
$ws.Range("H31").Value = 8932397
$ws.Range("I31").Value = 1998.5714
$ws.Range("K31").Value = 1998.5714
$ws.Range("M31").Value = -1703.5714
$ws.Range("H34").Value = 8932397
$ws.Range("I34").Value = 1998.5714
$ws.Range("K34").Value = 1998.5714
$ws.Range("M34").Value = -1796.5714
$ws.Range("H58").Value = 1567.1666
$ws.Range("I58").Value = 951.1667
$ws.Range("K58").Value = 951.1667
$ws.Range("M58").Value = -748.1667
$ws.Range("H99").Value = 3844.4666
$ws.Range("I99").Value = 3199.2
$ws.Range("J99").Value = 4167.1
$ws.Range("K99").Value = 3199.2
$ws.Range("L99").Value = 4167.1
$ws.Range("M99").Value = -1701.2
$ws.Range("N99").Value = -7163.1
$ws.Range("H113").Value = 748.0833
$ws.Range("I113").Value = 748.0833
$ws.Range("K113").Value = 748.0833
$ws.Range("M113").Value = 1421.9167
$ws.Range("H122").Value = 2948.4707
$ws.Range("I122").Value = 2090
$ws.Range("J122").Value = 4522.3335
$ws.Range("K122").Value = 6270
$ws.Range("L122").Value = 13567.0005
$ws.Range("M122").Value = -3820
$ws.Range("N122").Value = -18467.0005
$ws.Range("H126").Value = 3844.4666
$ws.Range("I126").Value = 3199.2
$ws.Range("J126").Value = 4167.1
$ws.Range("K126").Value = 9597.599999999999
$ws.Range("L126").Value = 12501.3
$ws.Range("M126").Value = -7127.599999999999
$ws.Range("N126").Value = -17441.3
$ws.Range("H132").Value = 3872.7878
$ws.Range("I132").Value = 3072.4827
$ws.Range("J132").Value = 9675
$ws.Range("K132").Value = 9217.4481
$ws.Range("L132").Value = 29025
$ws.Range("M132").Value = -6687.4481
$ws.Range("N132").Value = -34085
$ws.Range("H134").Value = 3698.7878
$ws.Range("I134").Value = 3588
$ws.Range("K134").Value = 10764
$ws.Range("M134").Value = -8229
$ws.Range("H136").Value = 1567.1666
$ws.Range("I136").Value = 951.1667
$ws.Range("K136").Value = 2853.5001
$ws.Range("M136").Value = -303.5001000000002

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 785.5
$ws.Range("J5").Value = 1237.5
$ws.Range("L5").Value = 3712.5
$ws.Range("N5").Value = -3936.5
$ws.Range("H131").Value = 3107456.2
$ws.Range("I131").Value = 22508.445
$ws.Range("J131").Value = 5883909.5
$ws.Range("K131").Value = 67525.33499999999
$ws.Range("L131").Value = 17651728.5
$ws.Range("M131").Value = -62485.33499999999
$ws.Range("N131").Value = -17661808.5
$ws.Range("H135").Value = 785.5
$ws.Range("J135").Value = 1237.5
$ws.Range("L135").Value = 11137.5
$ws.Range("N135").Value = -16207.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2067.0833
$ws.Range("I132").Value = 2010.0454
$ws.Range("K132").Value = 6030.1362
$ws.Range("M132").Value = -3500.1362

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4778.625
$ws.Range("I7").Value = 3545.8
$ws.Range("J7").Value = 6833.3335
$ws.Range("K7").Value = 3545.8
$ws.Range("L7").Value = 6833.3335
$ws.Range("M7").Value = -3433.8
$ws.Range("N7").Value = -7057.3335
$ws.Range("H122").Value = 3073
$ws.Range("I122").Value = 2226.5
$ws.Range("J122").Value = 4282.2856
$ws.Range("K122").Value = 6679.5
$ws.Range("L122").Value = 12846.8568
$ws.Range("M122").Value = -4229.5
$ws.Range("N122").Value = -17746.8568
$ws.Range("H126").Value = 4778.625
$ws.Range("I126").Value = 3545.8
$ws.Range("J126").Value = 6833.3335
$ws.Range("K126").Value = 10637.4
$ws.Range("L126").Value = 20500.0005
$ws.Range("M126").Value = -8167.400000000001
$ws.Range("N126").Value = -25440.0005
$ws.Range("H132").Value = 6397.4
$ws.Range("I132").Value = 7133
$ws.Range("K132").Value = 21399
$ws.Range("M132").Value = -18869

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2480.4546
$ws.Range("I136").Value = 2635.1924
$ws.Range("K136").Value = 7905.5772
$ws.Range("M136").Value = -5355.5772
